$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.482396960258484
$ws.Range("B1").Value = 3.604615211486816
$ws.Range("C1").Value = 2.965796232223511
$ws.Range("D1").Value = 1.383156180381775
$ws.Range("E1").Value = 0.7736344337463379
